$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revise the last two months of 2020 (rows 276 and 277)
$ws.Range("C276").Value = 214.461669921875
$ws.Range("C277").Value = 207.60650634765625

# Fill in 2021 data starting at row 278 (previously a blank trailing row)
$ws.Range("A278").Value = "'2021"
$ws.Range("B278").Value = 1
$ws.Range("C278").Value = 112.63906097412109

$ws.Range("A279").Value = "'2021"
$ws.Range("B279").Value = 2
$ws.Range("C279").Value = 138.00994873046875

$ws.Range("A280").Value = "'2021"
$ws.Range("B280").Value = 3
$ws.Range("C280").Value = 141.61357116699219

# New trailing row holding only the source-attribution note (mirrors the
# previous blank row 278 that used to sit after the last data row)
$ws.Range("A281").Value = "Source: 'Measuring Economic Policy Uncertainty' by Scott Baker, Nicholas Bloom and Steven J. Davis at www.PolicyUncertainty.com.  These data can be used freely with attribution to the authors, the paper, and the website."
$ws.Range("B281").Style = "Normal"
$ws.Range("C281").Style = "Normal"

# Clear the quote-prefix formatting Excel applies to the apostrophe-led
# text entries above so the cells fall back to the default (unstyled) cell
$ws.Range("A278:A280").Style = "Normal"
